# Add country partitioning to the narrative synthesis plot
#  - adds a "plot" boolean column (C) to Feuil1, pushing label/description to D/E
#  - rewrites the Kyoto Protocol row to the 1997-signing event (styled) and adds
#    several new timeline rows
#  - adds a second, blank sheet "Feuil2"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# 1. Plain (non shared-style) cell values: columns A, B, D, E for every row.
# ---------------------------------------------------------------------------

$rows = @(
    @{ n=1;  a='year_min'; b='year_max'; d='label';                                    e='description' },
    @{ n=2;  a=2015;       b=2015;       d='Paris Agreement';                          e='Aims to limit global warming to below 2 degrees Celsius above pre-industrial levels, with a goal of achieving a 1.5-degree ceiling. It encourages nations to strive for a 1.5-degree ceiling and requires substantial reductions in greenhouse gas emissions. ' },
    @{ n=3;  a=1997;       b=1997;       d='Kyoto Protocol';                           e='Signing of the Kyoto Protocol: the first international treaty to set legally binding obligations on developed countries to reduce greenhouse gas emissions (eneterd into force in 2005)' },
    @{ n=4;  a=1992;       b=1992;       d='UNFCCC';                                   e='United Nations Framework Convention on Climate Change' },
    @{ n=5;  a=2019;       b=2019;       d='European Green Deal';                      e='European Union''s ambitious roadmap to achieve climate neutrality by 2050' },
    @{ n=6;  a=2020;       b=2020;       d='COVID-19 pandemic';                        e='The economic shock instigated by the pandemic prompted regional and national policies to stimulate the economy and also deliver on climate mitigation commitments' },
    @{ n=7;  a=2014;       b=2014;       d='European 2030 climate & energy framework'; e='EU policy on energy economy, targetting reducing greenhouse gas emissions and increasing use of renewable energies.' },
    @{ n=8;  a=2006;       b=2006;       d='CDM';                                      e='The Clean Development Mechanism marked the uptick in carbon offsets (United Nations, 2018, Achievements of the Clean Development Mechanism, harnessing incentive for climate action (2001-2018))' },
    @{ n=9;  a=2008;       b=2009;       d='2008/9 Financial crisis';                  e='Financial crisis prompting green economic stimulus packages' },
    @{ n=10; a=2010;       b=2013;       d='IMO efficiency policies';                  e=$null },
    @{ n=11; a=2018;       b=2019;       d='IPCC Special Reports';                     e=$null }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.n, 1).Value = $r.a
    $ws.Cells.Item($r.n, 2).Value = $r.b
    $ws.Cells.Item($r.n, 4).Value = $r.d
    if ($null -ne $r.e) {
        $ws.Cells.Item($r.n, 5).Value = $r.e
    }
}

# Fix up the curly apostrophe / en-dash that plain ASCII source can't carry
# reliably through every console code path - set them precisely via .Characters.
$ws.Cells.Item(5, 5).Value = [char]0x0045 + 'uropean Union' + [char]0x2019 + 's ambitious roadmap to achieve climate neutrality by 2050'
$ws.Cells.Item(8, 5).Value = 'The Clean Development Mechanism marked the uptick in carbon offsets (United Nations, 2018, Achievements of the Clean Development Mechanism, harnessing incentive for climate action (2001' + [char]0x2013 + '2018))'

# ---------------------------------------------------------------------------
# 2. "plot" column (C): header text + TRUE/FALSE boolean-looking *text* values.
#    A literal Value = "TRUE" is auto-coerced to a real boolean by the engine
#    (same as Excel), so instead we write a text formula and freeze it to a
#    static value via copy / paste-special values - that keeps the shared
#    string as plain "TRUE"/"FALSE" text, exactly like the source file.
# ---------------------------------------------------------------------------

$ws.Cells.Item(1, 3).Value = "plot"

$plotVals = @{
    2 = "TRUE"; 3 = "FALSE"; 4 = "FALSE"; 5 = "FALSE"; 6 = "TRUE";
    7 = "FALSE"; 8 = "FALSE"; 9 = "TRUE"; 10 = "TRUE"; 11 = "TRUE"
}
foreach ($rn in $plotVals.Keys) {
    $ws.Cells.Item($rn, 3).Formula = '="' + $plotVals[$rn] + '"'
}
$ws.Range("C2:C11").Copy()
$ws.Range("C2:C11").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Styling
#    - A3:C3 (the re-dated Kyoto Protocol signing row) uses a dedicated font.
#    - E7 (European 2030 framework description) uses another dedicated font.
# ---------------------------------------------------------------------------

$ws.Range("E7").Font.Color = 2624516      # BGR(4,12,40)   -> RGB hex FF040C28

$ws.Range("A3:C3").Font.Name = "Arial"
$ws.Range("A3:C3").Font.Size = 10
$ws.Range("A3:C3").Font.Color = 2039583   # BGR(31,31,31)  -> RGB hex FF1F1F1F

# ---------------------------------------------------------------------------
# 4. Selection cursor, to match the post-edit state recorded in the workbook.
# ---------------------------------------------------------------------------

$ws.Range("D5").Select()

# ---------------------------------------------------------------------------
# 5. Add the second, empty worksheet "Feuil2" after Feuil1.
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Feuil2"

$wb.Worksheets.Item("Feuil1").Activate()
